$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "583. Delete Operation for Two Strings"
$ws.Range("B32").Value = "Medium"
$ws.Range("C32").Value = "https://leetcode.com/problems/delete-operation-for-two-strings/"
$ws.Range("D32").Value = 44510
$ws.Range("E32").Value = "二维动态规划"
$ws.Range("F32").Value = "看出是要求最长公共子序列就很常规了"
$ws.Range("G32").Value = "未复习"
